$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4733849228047293
$ws.Range("C2").Value = 0.620988910059761
$ws.Range("B3").Value = 2.358141528579331
$ws.Range("C3").Value = 1.39995674838293
$ws.Range("B4").Value = 3.746560921531258
$ws.Range("C4").Value = 2.336307608769147
$ws.Range("B5").Value = 11.59407588302485
$ws.Range("C5").Value = 3.04748142210458
$ws.Range("B6").Value = 12.86950713831006
$ws.Range("C6").Value = 3.799986921521622
$ws.Range("B7").Value = 13.00052058022762
$ws.Range("C7").Value = 4.702811255233944
$ws.Range("B8").Value = 23.16837075043347
$ws.Range("C8").Value = 5.457175793925451
$ws.Range("B9").Value = 23.33131498076093
$ws.Range("C9").Value = 6.26274671579596
$ws.Range("B10").Value = 26.35890585650083
$ws.Range("C10").Value = 7.090844081829997
$ws.Range("B11").Value = 27.56753992672493
$ws.Range("C11").Value = 7.822934826053727
$ws.Range("B12").Value = 37.27215380193955
$ws.Range("C12").Value = 8.597402731794086
$ws.Range("B13").Value = 37.40367343336951
$ws.Range("C13").Value = 9.449184605905581
$ws.Range("B14").Value = 39.77859304251324
$ws.Range("C14").Value = 9.989053787115578
$ws.Range("B15").Value = 41.38806295225672
$ws.Range("C15").Value = 10.93816672153917
$ws.Range("B16").Value = 45.9953815827279
$ws.Range("C16").Value = 11.86461632877548
$ws.Range("B17").Value = 46.10473026860484
$ws.Range("C17").Value = 12.67810243457209
$ws.Range("B18").Value = 50.06956776460472
$ws.Range("C18").Value = 13.55509340473673
$ws.Range("B19").Value = 50.16546488718263
$ws.Range("C19").Value = 14.26760160741861
$ws.Range("B20").Value = 54.47365246349317
$ws.Range("C20").Value = 15.25608698137854
$ws.Range("B21").Value = 54.58353742948654
$ws.Range("C21").Value = 16.02132569791805
$ws.Range("B22").Value = 55.40627128764899
$ws.Range("C22").Value = 16.85455078169796
$ws.Range("B23").Value = 55.49558018726816
$ws.Range("C23").Value = 17.52642440446178
$ws.Range("B24").Value = 58.3458078721814
$ws.Range("C24").Value = 18.03615387578817
$ws.Range("B25").Value = 58.5298218800695
$ws.Range("C25").Value = 18.93085835152302
$ws.Range("B26").Value = 59.87726225404547
$ws.Range("C26").Value = 19.43587576473157
$ws.Range("B27").Value = 60.49691183635471
$ws.Range("C27").Value = 20.76664317653444
$ws.Range("B28").Value = 60.5715289547451
$ws.Range("C28").Value = 21.4069081157275
$ws.Range("B29").Value = 63.32388201568553
$ws.Range("C29").Value = 22.4126073813558
$ws.Range("B30").Value = 66.73605128971755
$ws.Range("C30").Value = 23.32402776704072
$ws.Range("B31").Value = 67.78121761882453
$ws.Range("C31").Value = 23.96109288885926
$ws.Range("B32").Value = 68.01731728543184
$ws.Range("C32").Value = 24.78066817435118
$ws.Range("B33").Value = 68.15057952205673
$ws.Range("C33").Value = 25.57147229083767
$ws.Range("B34").Value = 70.72027583310772
$ws.Range("C34").Value = 26.31187928544972
$ws.Range("B35").Value = 70.81977578461461
$ws.Range("C35").Value = 27.06695510985734
$ws.Range("B36").Value = 71.52305812729104
$ws.Range("C36").Value = 27.74720835503489
$ws.Range("B37").Value = 78.11192235844234
$ws.Range("C37").Value = 28.76338176754271
$ws.Range("B38").Value = 78.22410253066282
$ws.Range("C38").Value = 29.50318046388015
$ws.Range("B39").Value = 79.53187560539477
$ws.Range("C39").Value = 30.51981960033121
$ws.Range("B40").Value = 79.71611421259716
$ws.Range("C40").Value = 31.42654254018872
$ws.Range("B41").Value = 81.14921628671817
$ws.Range("C41").Value = 32.44601523547607
$ws.Range("B42").Value = 84.81887098196647
$ws.Range("C42").Value = 33.47723609873309
$ws.Range("B43").Value = 88.30965237004179
$ws.Range("C43").Value = 34.32468299055528
$ws.Range("B44").Value = 88.90442744262754
$ws.Range("C44").Value = 34.90383199887786
$ws.Range("B45").Value = 89.26248886392386
$ws.Range("C45").Value = 36.05303214862865
$ws.Range("B46").Value = 96.32378327531636
$ws.Range("C46").Value = 36.5458970903201
$ws.Range("B47").Value = 96.4699439043366
$ws.Range("C47").Value = 37.3444260571822
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = 99.49515176249064
$ws.Range("C48").Value = 38.46059949670585

# Copy formatting (style) from A2 to the newly added A48 cell so it matches the
# bold/border/center-top style used by the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$excel.CutCopyMode = 0
